# Fruta / hortaliza, semanal
# Two new weekly data rows are inserted into the "Espárragos" data table:
#   - a new row at sheet row 3 (pushing the previous rows 3-26 down to 4-27)
#   - a new row at sheet row 11, i.e. after the insert above (pushing rows
#     that are now at 11-27 down to 12-28)
# The table grows from A1:R26 to A1:R28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 3 --------------------------------------------------
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44847
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1400
$ws.Range("N3").Value = "`$/kilo"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 1400
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"

# --- Insert new row 11 --------------------------------------------------
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 44848
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = "Espárragos"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 1300
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = 1400
$ws.Range("N11").Value = "`$/kilo"
$ws.Range("O11").Value = "Provincia de Diguillín"
$ws.Range("P11").Value = 1400
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"

# --- Fix up the Variedad swap between (now) rows 19 and 21 -------------
# Row 19 ("Verde") and row 21 ("Sin especificar") swap Variedad values
# compared to before the insert.
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("H21").Value = "Verde"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
